$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-01-09 14:20:15"

for ($row = 2; $row -le 514; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
